# Updated symbol list on Wed Feb 15 22:46:22 UTC 2023 with GitHub Actions
# Refreshes Price (D) and Volume(1h) (E) columns for the cryptos table.
# Values are written with a leading apostrophe to keep them as plain text
# (matching the workbook's existing inlineStr-as-text convention), and the
# style is reset to "Normal" afterwards so no numeric/text display format
# gets stamped onto the cell (Excel otherwise auto-applies a Text format
# when a numeric-looking string is force-entered as text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'314.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'6.20%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'44.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'7.12%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.142"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.95%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08061"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'6.57%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.520"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'2.76%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.683"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'5.73%"
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'17.09%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1295"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'6.71%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1917"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'4.40%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09396"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'4.52%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04250"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'5.87%"
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'-1.01%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001320"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'2.58%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005930"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.19%"
$ws.Range("E15").Style = "Normal"
$ws.Range("E17").Value = "'0.85%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'0.31%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3390"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'2.12%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'8.288"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'4.50%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1381"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-2.78%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.3138"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'4.52%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04222"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'3.97%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001277"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.81%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004552"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'14.95%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001344"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'9.14%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02696"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'11.89%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05441"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'4.34%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.005564"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-11.08%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007733"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.73%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1420"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'6.64%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007322"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-2.98%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008571"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'9.20%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3139"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'5.55%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006790"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'0.07%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000745"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.70%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'34.84%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.003976"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-5.43%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002087"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.70%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0001988"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.70%"
$ws.Range("E51").Style = "Normal"
